# Apply "Add data for 2022-09-04" update:
# - Rename sheet / update header label & title from 2022-08-26 to 2022-08-27
# - Update August (I9) and Total (I14) 2022 counts

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (updates <sheet name=.../> in workbook.xml)
$ws.Name = "Through 2022-08-27"

# Update the column header label (shared string used by I1)
$ws.Range("I1").Value = "2022 (through 08-27)"

# Update the data values
$ws.Range("I9").Value = 144
$ws.Range("I14").Value = 1115
